$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 53.7
$ws.Range("E3").Select()
